$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 becomes what used to be row 25's data, except B24 = 90843
$ws.Range("A24").Value = 112275263
$ws.Range("B24").Value = 90843
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 5448
$ws.Range("F24").Value = "Svartvit taggsvamp"
$ws.Range("G24").Value = "Phellodon connatus"
$ws.Range("H24").Value = "(Schultz) nom.prov"
$ws.Range("Q24").Value = 646200
$ws.Range("R24").Value = 6568598
$ws.Range("S24").Value = 50
$ws.Range("AL24").ClearContents()
$ws.Range("AO24").ClearContents()

# Row 25 becomes what used to be row 24's data, except B25 = 89749
$ws.Range("A25").Value = 112275250
$ws.Range("B25").Value = 89749
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 1106
$ws.Range("F25").Value = "Vågticka"
$ws.Range("G25").Value = "Osteina undosa"
$ws.Range("H25").Value = "(Peck) Zmitr."
$ws.Range("Q25").Value = 646166
$ws.Range("R25").Value = 6568529
$ws.Range("S25").Value = 25
$ws.Range("AL25").Value = "gran"
$ws.Range("AO25").Value = "gran"
